# Append the 2025-04-21 price row (row 51) to every price-history sheet in
# the Solar_Prices workbook, matching the "Updated Argent prices" commit.
#
# Each sheet already stores Date (col A) / Price (col B) as text values for
# every existing row, so the new row has to be written as text too (not an
# Excel date-serial / number) to stay consistent with the rest of the column.

$wb = $excel.ActiveWorkbook

$newDate = "2025-04-21"

# Sheet name -> new Price value for row 51 (col B), as shown in the diff.
$updates = [ordered]@{
    "N-Dense"                   = "39.5"
    "N-Type"                    = "40"
    "N-type Wafer"               = "1.23"
    "Cell Topcon 183mm"          = "0.293"
    "Module Topcon 183mm"        = "0.09"
    "Silver Rear_side"           = "5,329"
    "Silver Busbar front-side"   = "7,977"
    "Silver finger front-side"   = "8,027"
    "USD_CNY"                    = "7.3173"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $price = $updates[$sheetName]

    $dateCell = $ws.Cells.Item(51, 1)
    $priceCell = $ws.Cells.Item(51, 2)

    # Force text storage (otherwise "2025-04-21" / "39.5" etc. get parsed as
    # a real date serial / number by Excel's input parser), then drop back to
    # the default "Normal" style so no stray number-format style is left on
    # the new cells.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDate
    $dateCell.Style = "Normal"

    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"
}
